$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 37605
$ws.Range("D2").Value = 54390950
$ws.Range("C3").Value = 90722
$ws.Range("D3").Value = 132997555
$ws.Range("C4").Value = 31103
$ws.Range("D4").Value = 46063454
$ws.Range("C5").Value = 8673
$ws.Range("D5").Value = 12891203
$ws.Range("C11").Value = 41171
$ws.Range("D11").Value = 55872429
$ws.Range("C12").Value = 9627
$ws.Range("D12").Value = 13925099
$ws.Range("C13").Value = 25885
$ws.Range("D13").Value = 37964455
$ws.Range("C14").Value = 8301
$ws.Range("D14").Value = 12319718
$ws.Range("C15").Value = 2146
$ws.Range("D15").Value = 3191383
$ws.Range("C19").Value = 10189
$ws.Range("D19").Value = 13496811
$ws.Range("C20").Value = 13349
$ws.Range("D20").Value = 19277099
$ws.Range("C21").Value = 31581
$ws.Range("D21").Value = 46348864
$ws.Range("C22").Value = 10207
$ws.Range("D22").Value = 15173555
$ws.Range("C23").Value = 2632
$ws.Range("D23").Value = 3913182
$ws.Range("C26").Value = 11651
$ws.Range("D26").Value = 15564923
$ws.Range("C27").Value = 7621
$ws.Range("D27").Value = 11041098
$ws.Range("C28").Value = 22430
$ws.Range("D28").Value = 32923931
$ws.Range("C29").Value = 7793
$ws.Range("D29").Value = 11597133
$ws.Range("C30").Value = 1954
$ws.Range("D30").Value = 2915499
$ws.Range("C33").Value = 8271
$ws.Range("D33").Value = 10930630
$ws.Range("C34").Value = 3229
$ws.Range("D34").Value = 4660637
$ws.Range("C35").Value = 7803
$ws.Range("D35").Value = 11395016
$ws.Range("C36").Value = 3168
$ws.Range("D36").Value = 4694961
$ws.Range("C38").Value = 163
$ws.Range("D38").Value = 242732
$ws.Range("C40").Value = 2457
$ws.Range("D40").Value = 3320816
$ws.Range("C41").Value = 17179
$ws.Range("D41").Value = 24844402
$ws.Range("C42").Value = 50955
$ws.Range("D42").Value = 74702656
$ws.Range("C43").Value = 18969
$ws.Range("D43").Value = 28175943
$ws.Range("C44").Value = 5592
$ws.Range("D44").Value = 8327978
$ws.Range("C45").Value = 1198
$ws.Range("D45").Value = 1787545
$ws.Range("C49").Value = 16625
$ws.Range("D49").Value = 22143593
$ws.Range("C50").Value = 1997
$ws.Range("D50").Value = 2897562
$ws.Range("C51").Value = 6835
$ws.Range("D51").Value = 10048699
$ws.Range("C52").Value = 2335
$ws.Range("D52").Value = 3487418
$ws.Range("C53").Value = 750
$ws.Range("D53").Value = 1120305
$ws.Range("C54").Value = 184
$ws.Range("D54").Value = 272833
$ws.Range("C56").Value = 6817
$ws.Range("D56").Value = 9386508
$ws.Range("C57").Value = 928
$ws.Range("D57").Value = 1361579
$ws.Range("C58").Value = 2328
$ws.Range("D58").Value = 3450917
$ws.Range("C59").Value = 930
$ws.Range("D59").Value = 1384501
$ws.Range("C60").Value = 318
$ws.Range("D60").Value = 476758
$ws.Range("C61").Value = 101
$ws.Range("D61").Value = 151350
$ws.Range("C63").Value = 1370
$ws.Range("D63").Value = 1926206
$ws.Range("C64").Value = 15296
$ws.Range("D64").Value = 22097128
$ws.Range("C65").Value = 44567
$ws.Range("D65").Value = 65221783
$ws.Range("C66").Value = 15663
$ws.Range("D66").Value = 23278932
$ws.Range("C67").Value = 4559
$ws.Range("D67").Value = 6790792
$ws.Range("C68").Value = 916
$ws.Range("D68").Value = 1362168
$ws.Range("C72").Value = 15040
$ws.Range("D72").Value = 19836353
$ws.Range("C73").Value = 51085
$ws.Range("D73").Value = 74342597
$ws.Range("C74").Value = 145305
$ws.Range("D74").Value = 214075977
$ws.Range("C75").Value = 63368
$ws.Range("D75").Value = 94427770
$ws.Range("C76").Value = 20234
$ws.Range("D76").Value = 30231709
$ws.Range("C77").Value = 4785
$ws.Range("D77").Value = 7149223
$ws.Range("C78").Value = 263
$ws.Range("D78").Value = 389670
$ws.Range("C84").Value = 50532
$ws.Range("D84").Value = 68760103
$ws.Range("C85").Value = 4564
$ws.Range("D85").Value = 6612669
$ws.Range("C86").Value = 11525
$ws.Range("D86").Value = 16932440
$ws.Range("C87").Value = 3865
$ws.Range("D87").Value = 5760083
$ws.Range("C88").Value = 1342
$ws.Range("D88").Value = 2005489
$ws.Range("C92").Value = 5370
$ws.Range("D92").Value = 7221585
$ws.Range("C93").Value = 1587
$ws.Range("D93").Value = 2285432
$ws.Range("C94").Value = 5121
$ws.Range("D94").Value = 7543139
$ws.Range("C95").Value = 1934
$ws.Range("D95").Value = 2880937
$ws.Range("C96").Value = 686
$ws.Range("D96").Value = 1027960
$ws.Range("C97").Value = 182
$ws.Range("D97").Value = 272113
$ws.Range("C100").Value = 3519
$ws.Range("D100").Value = 4658469
$ws.Range("C101").Value = 593
$ws.Range("D101").Value = 883164
$ws.Range("C102").Value = 348
$ws.Range("D102").Value = 519530
$ws.Range("C106").Value = 10723
$ws.Range("D106").Value = 15556887
$ws.Range("C107").Value = 29109
$ws.Range("D107").Value = 42771513
$ws.Range("C108").Value = 9757
$ws.Range("D108").Value = 14509150
$ws.Range("C110").Value = 488
$ws.Range("D110").Value = 727046
$ws.Range("C111").Value = 50
$ws.Range("D111").Value = 75000
$ws.Range("C113").Value = 9755
$ws.Range("D113").Value = 12889313
$ws.Range("C114").Value = 30309
$ws.Range("D114").Value = 43711283
$ws.Range("C115").Value = 65938
$ws.Range("D115").Value = 96508211
$ws.Range("C116").Value = 21305
$ws.Range("D116").Value = 31662332
$ws.Range("C117").Value = 6044
$ws.Range("D117").Value = 9005354
$ws.Range("C123").Value = 25745
$ws.Range("D123").Value = 34403307
$ws.Range("C124").Value = 35832
$ws.Range("D124").Value = 51722895
$ws.Range("C125").Value = 76529
$ws.Range("D125").Value = 111916050
$ws.Range("C126").Value = 23774
$ws.Range("D126").Value = 35286600
$ws.Range("C127").Value = 6370
$ws.Range("D127").Value = 9466051
$ws.Range("C128").Value = 1227
$ws.Range("D128").Value = 1824911
$ws.Range("C132").Value = 31650
$ws.Range("D132").Value = 42046604
$ws.Range("C133").Value = 13174
$ws.Range("D133").Value = 19071312
$ws.Range("C134").Value = 32216
$ws.Range("D134").Value = 47322428
$ws.Range("C135").Value = 11450
$ws.Range("D135").Value = 17013042
$ws.Range("C136").Value = 2949
$ws.Range("D136").Value = 4396305
$ws.Range("C140").Value = 10776
$ws.Range("D140").Value = 14373804
$ws.Range("C141").Value = 34884
$ws.Range("D141").Value = 50379267
$ws.Range("C142").Value = 80951
$ws.Range("D142").Value = 118613156
$ws.Range("C143").Value = 24274
$ws.Range("D143").Value = 36065555
$ws.Range("C144").Value = 6366
$ws.Range("D144").Value = 9498567
$ws.Range("C148").Value = 29048
$ws.Range("D148").Value = 39204787
